# Data update 06/10/2024: the sheet only needs to keep the "link" header
# in A1 - all the previously collected Facebook post URLs (rows 2:65) are
# removed from the sheet entirely (not just cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

if ($lastRow -gt 1) {
    $ws.Range("A2:A$lastRow").EntireRow.Delete()
}

# Reset the window's scroll position / selection the way Excel would after
# trimming the sheet back down to a single row.
[void]$ws.Range("A5").Select()

# Match the saved window size recorded the next time the workbook was opened.
$win = $excel.ActiveWindow
$win.Width = 27945
$win.Height = 12180
